$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D25").Value = "\n<Lime>I knew it!"
$ws.Range("D26").Value = "\n<Lime>I was just thinking, `"there's was no time to escape at all-`".`nYou hid somewhere, right?So I thought to wait here, just in case!"
$ws.Range("D27").Value = "\n<Lime>Ehehehehe-♥`nCaught you-♥"
$ws.Range("D28").Value = "\n<\n[3]>It's hard to move when everything is so slippery right?`nIf you don't escape soon you'll let out your white stuff you know?"
$ws.Range("D29").Value = "\C[3]※A timing bar will display if you're caught.`n\C[0]Go ahead and stop in the red or yellow areas.Red areas will reduce damage taken by half."
$ws.Range("D30").Value = "\n<\n[3]>Squish-...♥`nThey're so soft you're about to melt aren't you?♥Feels so good- Feels so good-♥"
$ws.Range("D31").Value = "\n<\n[3]>What's this? You're about to cum already-?`nSo you couldn't resist my boobs after all? Ahaha♥Well then, I'll give you one more pump, let it all out♪"
$ws.Range("D32").Value = "\n<\n[3]>Ahh ahh♥ It's shooting between my boobies-♥`nYour penis seems delighted too♥ I'm so happy-♥"
$ws.Range("D34").Value = "\n<Lime>Whoa- Were you seriously not trying to resist at all-?`nAnd the game has only just started, next time I'llmilk you for real okay?"
$ws.Range("D35").Value = "\n<Lime>Ah! You escaped...`nOh well-."
$ws.Range("D36").Value = "\n<Lime>That is...you heard what Lily said right?`nAbout this escape game thingy..."
$ws.Range("D37").Value = "\n<Lime>This mansion is very, very big-.`nI don't think you'll get out...but feel free to try your hardest,and search everywhere you can okay?"
$ws.Range("D38").Value = "\n<Lime>Lily and Shina seem to be happy after all.`nSo you shouldn't give up so easily OK? Enjoy yourselfas much as possible now-"
$ws.Range("D39").Value = "\n<Lime>Next time you cum, I won't go easy on you now alright?`nTeehee♥ Bye bye♥"
$ws.Range("D40").Value = "\n<Lime>Now then, do your best okay?`nNext time I'll make you go pew pew for realsies-"
